$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 51: apply data from old row 52 (permutation within group)
$ws.Range("B51").Value = 5140743
$ws.Range("F51").Value = "Stal Rzeszow"
$ws.Range("G51").Value = "Termalica BB Nieciecza"
$ws.Range("H51").Value = 2
$ws.Range("I51").Value = 2
$ws.Range("J51").Value = "D"
$ws.Range("K51").Value = 3
$ws.Range("L51").Value = 3.3
$ws.Range("M51").Value = 2.2
$ws.Range("N51").Value = 2.9
$ws.Range("O51").Value = 3.3
$ws.Range("P51").Value = 2.25
$ws.Range("Q51").Value = 0.25
$ws.Range("R51").Value = 1.825
$ws.Range("S51").Value = 1.975
$ws.Range("T51").Value = 2.5
$ws.Range("U51").Value = 1.95
$ws.Range("V51").Value = 1.85
$ws.Range("W51").Value = -1
$ws.Range("X51").Value = 2.3
$ws.Range("Y51").Value = -1
$ws.Range("Z51").Value = 0.4125
$ws.Range("AA51").Value = -0.5
$ws.Range("AB51").Value = 0.95
$ws.Range("AC51").Value = -1

# Row 52: apply data from old row 54 (permutation within group)
$ws.Range("B52").Value = 5139053
$ws.Range("F52").Value = "Chrobry Glogow"
$ws.Range("G52").Value = "Zaglebie Sosnowiec"
$ws.Range("H52").Value = 0
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = "D"
$ws.Range("K52").Value = 2.45
$ws.Range("L52").Value = 3.2
$ws.Range("M52").Value = 2.55
$ws.Range("N52").Value = 2.7
$ws.Range("O52").Value = 3.2
$ws.Range("P52").Value = 2.375
$ws.Range("Q52").Value = 0
$ws.Range("R52").Value = 2.05
$ws.Range("S52").Value = 1.75
$ws.Range("T52").Value = 2.25
$ws.Range("U52").Value = 1.875
$ws.Range("V52").Value = 1.925
$ws.Range("W52").Value = -1
$ws.Range("X52").Value = 2.2
$ws.Range("Y52").Value = -1
$ws.Range("Z52").Value = 0
$ws.Range("AA52").Value = -0
$ws.Range("AB52").Value = -1
$ws.Range("AC52").Value = 0.925

# Row 53: apply data from old row 51 (permutation within group)
$ws.Range("B53").Value = 5139054
$ws.Range("F53").Value = "GKS Tychy 71"
$ws.Range("G53").Value = "Sandecja Nowy Sacz"
$ws.Range("H53").Value = 2
$ws.Range("I53").Value = 3
$ws.Range("J53").Value = "A"
$ws.Range("K53").Value = 2.15
$ws.Range("L53").Value = 3.2
$ws.Range("M53").Value = 3.1
$ws.Range("N53").Value = 2.375
$ws.Range("O53").Value = 3
$ws.Range("P53").Value = 3
$ws.Range("Q53").Value = -0.25
$ws.Range("R53").Value = 2.025
$ws.Range("S53").Value = 1.775
$ws.Range("T53").Value = 2.25
$ws.Range("U53").Value = 1.975
$ws.Range("V53").Value = 1.825
$ws.Range("W53").Value = -1
$ws.Range("X53").Value = -1
$ws.Range("Y53").Value = 2
$ws.Range("Z53").Value = -1
$ws.Range("AA53").Value = 0.7749999999999999
$ws.Range("AB53").Value = 0.9750000000000001
$ws.Range("AC53").Value = -1

# Row 54: apply data from old row 53 (permutation within group)
$ws.Range("B54").Value = 5139056
$ws.Range("F54").Value = "Odra Opole"
$ws.Range("G54").Value = "Arka Gdynia"
$ws.Range("H54").Value = 1
$ws.Range("I54").Value = 1
$ws.Range("J54").Value = "D"
$ws.Range("K54").Value = 3.75
$ws.Range("L54").Value = 3.5
$ws.Range("M54").Value = 1.85
$ws.Range("N54").Value = 3.4
$ws.Range("O54").Value = 3.5
$ws.Range("P54").Value = 1.909
$ws.Range("Q54").Value = 0.5
$ws.Range("R54").Value = 1.85
$ws.Range("S54").Value = 2
$ws.Range("T54").Value = 2.75
$ws.Range("U54").Value = 2
$ws.Range("V54").Value = 1.85
$ws.Range("W54").Value = -1
$ws.Range("X54").Value = 2.5
$ws.Range("Y54").Value = -1
$ws.Range("Z54").Value = 0.8500000000000001
$ws.Range("AA54").Value = -1
$ws.Range("AB54").Value = -1
$ws.Range("AC54").Value = 0.8500000000000001

# Row 136: apply data from old row 138 (permutation within group)
$ws.Range("B136").Value = 5448048
$ws.Range("F136").Value = "Zaglebie Sosnowiec"
$ws.Range("G136").Value = "Sandecja Nowy Sacz"
$ws.Range("H136").Value = 1
$ws.Range("I136").Value = 1
$ws.Range("J136").Value = "D"
$ws.Range("K136").Value = 2.1
$ws.Range("L136").Value = 3.2
$ws.Range("M136").Value = 3.3
$ws.Range("N136").Value = 2.1
$ws.Range("O136").Value = 3.2
$ws.Range("P136").Value = 3.1
$ws.Range("Q136").Value = -0.25
$ws.Range("R136").Value = 1.875
$ws.Range("S136").Value = 1.925
$ws.Range("T136").Value = 2.25
$ws.Range("U136").Value = 1.85
$ws.Range("V136").Value = 1.95
$ws.Range("W136").Value = -1
$ws.Range("X136").Value = 2.2
$ws.Range("Y136").Value = -1
$ws.Range("Z136").Value = -0.5
$ws.Range("AA136").Value = 0.4625
$ws.Range("AB136").Value = -0.5
$ws.Range("AC136").Value = 0.475

# Row 137: apply data from old row 136 (permutation within group)
$ws.Range("B137").Value = 5451608
$ws.Range("F137").Value = "Termalica BB Nieciecza"
$ws.Range("G137").Value = "Arka Gdynia"
$ws.Range("H137").Value = 2
$ws.Range("I137").Value = 1
$ws.Range("J137").Value = "H"
$ws.Range("K137").Value = 1.909
$ws.Range("L137").Value = 3.5
$ws.Range("M137").Value = 3.5
$ws.Range("N137").Value = 1.909
$ws.Range("O137").Value = 3.5
$ws.Range("P137").Value = 3.6
$ws.Range("Q137").Value = -0.5
$ws.Range("R137").Value = 1.95
$ws.Range("S137").Value = 1.85
$ws.Range("T137").Value = 2.75
$ws.Range("U137").Value = 1.8
$ws.Range("V137").Value = 2
$ws.Range("W137").Value = 0.909
$ws.Range("X137").Value = -1
$ws.Range("Y137").Value = -1
$ws.Range("Z137").Value = 0.95
$ws.Range("AA137").Value = -1
$ws.Range("AB137").Value = 0.4
$ws.Range("AC137").Value = -0.5

# Row 138: apply data from old row 140 (permutation within group)
$ws.Range("B138").Value = 5447925
$ws.Range("F138").Value = "Gornik Leczna"
$ws.Range("G138").Value = "Wisla Krakow"
$ws.Range("H138").Value = 0
$ws.Range("I138").Value = 3
$ws.Range("J138").Value = "A"
$ws.Range("K138").Value = 5.5
$ws.Range("L138").Value = 4
$ws.Range("M138").Value = 1.5
$ws.Range("N138").Value = 4.5
$ws.Range("O138").Value = 4
$ws.Range("P138").Value = 1.615
$ws.Range("Q138").Value = 0.75
$ws.Range("R138").Value = 2.05
$ws.Range("S138").Value = 1.8
$ws.Range("T138").Value = 3
$ws.Range("U138").Value = 2
$ws.Range("V138").Value = 1.85
$ws.Range("W138").Value = -1
$ws.Range("X138").Value = -1
$ws.Range("Y138").Value = 0.615
$ws.Range("Z138").Value = -1
$ws.Range("AA138").Value = 0.8
$ws.Range("AB138").Value = 0
$ws.Range("AC138").Value = -0

# Row 139: apply data from old row 142 (permutation within group)
$ws.Range("B139").Value = 5452381
$ws.Range("F139").Value = "MKS Puszcza Niepolomice"
$ws.Range("G139").Value = "Chrobry Glogow"
$ws.Range("H139").Value = 0
$ws.Range("I139").Value = 1
$ws.Range("J139").Value = "A"
$ws.Range("K139").Value = 1.571
$ws.Range("L139").Value = 4
$ws.Range("M139").Value = 5
$ws.Range("N139").Value = 1.4
$ws.Range("O139").Value = 4.5
$ws.Range("P139").Value = 6.5
$ws.Range("Q139").Value = -1.25
$ws.Range("R139").Value = 1.9
$ws.Range("S139").Value = 1.95
$ws.Range("T139").Value = 3
$ws.Range("U139").Value = 2.025
$ws.Range("V139").Value = 1.825
$ws.Range("W139").Value = -1
$ws.Range("X139").Value = -1
$ws.Range("Y139").Value = 5.5
$ws.Range("Z139").Value = -1
$ws.Range("AA139").Value = 0.95
$ws.Range("AB139").Value = -1
$ws.Range("AC139").Value = 0.825

# Row 140: apply data from old row 143 (permutation within group)
$ws.Range("B140").Value = 5451609
$ws.Range("F140").Value = "Ruch Chorzow"
$ws.Range("G140").Value = "GKS Tychy 71"
$ws.Range("H140").Value = 1
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = "H"
$ws.Range("K140").Value = 1.3
$ws.Range("L140").Value = 5
$ws.Range("M140").Value = 7.5
$ws.Range("N140").Value = 1.333
$ws.Range("O140").Value = 4.75
$ws.Range("P140").Value = 8
$ws.Range("Q140").Value = -1.5
$ws.Range("R140").Value = 2
$ws.Range("S140").Value = 1.8
$ws.Range("T140").Value = 2.5
$ws.Range("U140").Value = 1.825
$ws.Range("V140").Value = 1.975
$ws.Range("W140").Value = 0.333
$ws.Range("X140").Value = -1
$ws.Range("Y140").Value = -1
$ws.Range("Z140").Value = -1
$ws.Range("AA140").Value = 0.8
$ws.Range("AB140").Value = -1
$ws.Range("AC140").Value = 0.9750000000000001

# Row 141: apply data from old row 137 (permutation within group)
$ws.Range("B141").Value = 5451607
$ws.Range("F141").Value = "Podbeskidzie Bielsko Biala"
$ws.Range("G141").Value = "Resovia Rzeszow"
$ws.Range("H141").Value = 4
$ws.Range("I141").Value = 3
$ws.Range("J141").Value = "H"
$ws.Range("K141").Value = 1.615
$ws.Range("L141").Value = 3.75
$ws.Range("M141").Value = 4.75
$ws.Range("N141").Value = 1.363
$ws.Range("O141").Value = 4.75
$ws.Range("P141").Value = 6
$ws.Range("Q141").Value = -1.25
$ws.Range("R141").Value = 1.825
$ws.Range("S141").Value = 1.975
$ws.Range("T141").Value = 3.25
$ws.Range("U141").Value = 1.925
$ws.Range("V141").Value = 1.875
$ws.Range("W141").Value = 0.363
$ws.Range("X141").Value = -1
$ws.Range("Y141").Value = -1
$ws.Range("Z141").Value = -0.5
$ws.Range("AA141").Value = 0.4875
$ws.Range("AB141").Value = 0.925
$ws.Range("AC141").Value = -1

# Row 142: apply data from old row 141 (permutation within group)
$ws.Range("B142").Value = 5451610
$ws.Range("F142").Value = "Stal Rzeszow"
$ws.Range("G142").Value = "Skra Czestochowa"
$ws.Range("H142").Value = 2
$ws.Range("I142").Value = 1
$ws.Range("J142").Value = "H"
$ws.Range("K142").Value = 1.444
$ws.Range("L142").Value = 4
$ws.Range("M142").Value = 6.5
$ws.Range("N142").Value = 1.333
$ws.Range("O142").Value = 4.333
$ws.Range("P142").Value = 8
$ws.Range("Q142").Value = -1.5
$ws.Range("R142").Value = 1.95
$ws.Range("S142").Value = 1.85
$ws.Range("T142").Value = 2.75
$ws.Range("U142").Value = 1.875
$ws.Range("V142").Value = 1.925
$ws.Range("W142").Value = 0.333
$ws.Range("X142").Value = -1
$ws.Range("Y142").Value = -1
$ws.Range("Z142").Value = -1
$ws.Range("AA142").Value = 0.8500000000000001
$ws.Range("AB142").Value = 0.4375
$ws.Range("AC142").Value = -0.5

# Row 143: apply data from old row 144 (permutation within group)
$ws.Range("B143").Value = 5448049
$ws.Range("F143").Value = "LKS Lodz"
$ws.Range("G143").Value = "Odra Opole"
$ws.Range("H143").Value = 1
$ws.Range("I143").Value = 0
$ws.Range("J143").Value = "H"
$ws.Range("K143").Value = 1.571
$ws.Range("L143").Value = 3.75
$ws.Range("M143").Value = 5
$ws.Range("N143").Value = 1.444
$ws.Range("O143").Value = 4
$ws.Range("P143").Value = 6
$ws.Range("Q143").Value = -1
$ws.Range("R143").Value = 1.775
$ws.Range("S143").Value = 2.025
$ws.Range("T143").Value = 2.75
$ws.Range("U143").Value = 1.9
$ws.Range("V143").Value = 1.9
$ws.Range("W143").Value = 0.444
$ws.Range("X143").Value = -1
$ws.Range("Y143").Value = -1
$ws.Range("Z143").Value = 0
$ws.Range("AA143").Value = -0
$ws.Range("AB143").Value = -1
$ws.Range("AC143").Value = 0.8999999999999999

# Row 144: apply data from old row 139 (permutation within group)
$ws.Range("B144").Value = 5448050
$ws.Range("F144").Value = "Chojniczanka Chojnice"
$ws.Range("G144").Value = "GKS Katowice"
$ws.Range("H144").Value = 3
$ws.Range("I144").Value = 3
$ws.Range("J144").Value = "D"
$ws.Range("K144").Value = 2.75
$ws.Range("L144").Value = 3.25
$ws.Range("M144").Value = 2.375
$ws.Range("N144").Value = 2.4
$ws.Range("O144").Value = 3.25
$ws.Range("P144").Value = 2.7
$ws.Range("Q144").Value = 0
$ws.Range("R144").Value = 1.8
$ws.Range("S144").Value = 2.05
$ws.Range("T144").Value = 2.5
$ws.Range("U144").Value = 2
$ws.Range("V144").Value = 1.85
$ws.Range("W144").Value = -1
$ws.Range("X144").Value = 2.25
$ws.Range("Y144").Value = -1
$ws.Range("Z144").Value = 0
$ws.Range("AA144").Value = -0
$ws.Range("AB144").Value = 1
$ws.Range("AC144").Value = -1

# Row 255: apply data from old row 256 (permutation within group)
$ws.Range("B255").Value = 6805700
$ws.Range("F255").Value = "GKS Katowice"
$ws.Range("G255").Value = "Polonia Warsaw"
$ws.Range("H255").Value = 0
$ws.Range("I255").Value = 2
$ws.Range("J255").Value = "A"
$ws.Range("K255").Value = 1.7
$ws.Range("L255").Value = 3.75
$ws.Range("M255").Value = 4.2
$ws.Range("N255").Value = 1.75
$ws.Range("O255").Value = 3.6
$ws.Range("P255").Value = 4
$ws.Range("Q255").Value = -0.5
$ws.Range("R255").Value = 1.775
$ws.Range("S255").Value = 2.025
$ws.Range("T255").Value = 2.75
$ws.Range("U255").Value = 2
$ws.Range("V255").Value = 1.8
$ws.Range("W255").Value = -1
$ws.Range("X255").Value = -1
$ws.Range("Y255").Value = 3
$ws.Range("Z255").Value = -1
$ws.Range("AA255").Value = 1.025
$ws.Range("AB255").Value = -1
$ws.Range("AC255").Value = 0.8

# Row 256: apply data from old row 255 (permutation within group)
$ws.Range("B256").Value = 6803761
$ws.Range("F256").Value = "Wisla Plock"
$ws.Range("G256").Value = "Chrobry Glogow"
$ws.Range("H256").Value = 2
$ws.Range("I256").Value = 1
$ws.Range("J256").Value = "H"
$ws.Range("K256").Value = 1.5
$ws.Range("L256").Value = 4
$ws.Range("M256").Value = 5.25
$ws.Range("N256").Value = 1.5
$ws.Range("O256").Value = 4
$ws.Range("P256").Value = 5.5
$ws.Range("Q256").Value = -1
$ws.Range("R256").Value = 1.9
$ws.Range("S256").Value = 1.9
$ws.Range("T256").Value = 2.75
$ws.Range("U256").Value = 1.975
$ws.Range("V256").Value = 1.825
$ws.Range("W256").Value = 0.5
$ws.Range("X256").Value = -1
$ws.Range("Y256").Value = -1
$ws.Range("Z256").Value = 0
$ws.Range("AA256").Value = -0
$ws.Range("AB256").Value = 0.4875
$ws.Range("AC256").Value = -0.5

# Row 302: apply data from old row 303 (permutation within group)
$ws.Range("B302").Value = 6803794
$ws.Range("F302").Value = "Wisla Krakow"
$ws.Range("G302").Value = "Gornik Leczna"
$ws.Range("H302").Value = 4
$ws.Range("I302").Value = 0
$ws.Range("J302").Value = "H"
$ws.Range("K302").Value = 1.4
$ws.Range("L302").Value = 4.75
$ws.Range("M302").Value = 7
$ws.Range("N302").Value = 1.363
$ws.Range("O302").Value = 4.75
$ws.Range("P302").Value = 7.5
$ws.Range("Q302").Value = -1.25
$ws.Range("R302").Value = 1.8
$ws.Range("S302").Value = 2
$ws.Range("T302").Value = 2.75
$ws.Range("U302").Value = 1.775
$ws.Range("V302").Value = 2.025
$ws.Range("W302").Value = 0.363
$ws.Range("X302").Value = -1
$ws.Range("Y302").Value = -1
$ws.Range("Z302").Value = 0.8
$ws.Range("AA302").Value = -1
$ws.Range("AB302").Value = 0.7749999999999999
$ws.Range("AC302").Value = -1

# Row 303: apply data from old row 302 (permutation within group)
$ws.Range("B303").Value = 6803793
$ws.Range("F303").Value = "Odra Opole"
$ws.Range("G303").Value = "Stal Rzeszow"
$ws.Range("H303").Value = 1
$ws.Range("I303").Value = 1
$ws.Range("J303").Value = "D"
$ws.Range("K303").Value = 2.05
$ws.Range("L303").Value = 3.4
$ws.Range("M303").Value = 3.5
$ws.Range("N303").Value = 2.1
$ws.Range("O303").Value = 3.4
$ws.Range("P303").Value = 3.4
$ws.Range("Q303").Value = -0.25
$ws.Range("R303").Value = 1.825
$ws.Range("S303").Value = 1.975
$ws.Range("T303").Value = 2.5
$ws.Range("U303").Value = 1.95
$ws.Range("V303").Value = 1.85
$ws.Range("W303").Value = -1
$ws.Range("X303").Value = 2.4
$ws.Range("Y303").Value = -1
$ws.Range("Z303").Value = -0.5
$ws.Range("AA303").Value = 0.4875
$ws.Range("AB303").Value = -1
$ws.Range("AC303").Value = 0.8500000000000001

# Row 316: fill in match result (previously unplayed fixture)
$ws.Range("H316").Value = 0
$ws.Range("I316").Value = 0
$ws.Range("J316").Value = "D"
$ws.Range("N316").Value = 4.333
$ws.Range("O316").Value = 3.25
$ws.Range("P316").Value = 1.909
$ws.Range("Q316").Value = 0.5
$ws.Range("R316").Value = 1.875
$ws.Range("S316").Value = 1.925
$ws.Range("U316").Value = 1.875
$ws.Range("V316").Value = 1.925
$ws.Range("W316").Value = -1
$ws.Range("X316").Value = 2.25
$ws.Range("Y316").Value = -1
$ws.Range("Z316").Value = 0.875
$ws.Range("AA316").Value = -1
$ws.Range("AB316").Value = -1
$ws.Range("AC316").Value = 0.925

# Row 317: fill in match result (previously unplayed fixture)
$ws.Range("H317").Value = 2
$ws.Range("I317").Value = 0
$ws.Range("J317").Value = "H"
$ws.Range("N317").Value = 2.2
$ws.Range("P317").Value = 3.3
$ws.Range("R317").Value = 1.825
$ws.Range("S317").Value = 1.975
$ws.Range("V317").Value = 1.825
$ws.Range("W317").Value = 1.2
$ws.Range("X317").Value = -1
$ws.Range("Y317").Value = -1
$ws.Range("Z317").Value = 0.825
$ws.Range("AA317").Value = -1
$ws.Range("AB317").Value = -0.5
$ws.Range("AC317").Value = 0.4125

# Row 320: odds update only
$ws.Range("N320").Value = 2.15
$ws.Range("P320").Value = 3.25
$ws.Range("R320").Value = 1.9
$ws.Range("S320").Value = 1.95

# Row 323: odds update only
$ws.Range("N323").Value = 2.45
$ws.Range("O323").Value = 3.1
$ws.Range("P323").Value = 2.9
$ws.Range("R323").Value = 2.125
$ws.Range("S323").Value = 1.75

# Row 324: odds update only
$ws.Range("N324").Value = 2.45
$ws.Range("P324").Value = 2.9
$ws.Range("R324").Value = 2.125
$ws.Range("S324").Value = 1.75
